$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "48.402.77"
$ws.Range("E2").Value2 = "  +2.57%  "

$ws.Range("D3").Value2 = "2.525.40"
$ws.Range("E3").Value2 = "  +1.61%  "

$ws.Range("E4").Value2 = "  +0.06%  "

$ws.Range("D5").Value2 = "110.29"
$ws.Range("E5").Value2 = "  +2.50%  "

$ws.Range("D6").Value2 = "322.53"
$ws.Range("E6").Value2 = "  +0.62%  "

$ws.Range("E7").Value2 = "  +2.34%  "

$ws.Range("E8").Value2 = "  +0.06%  "

$ws.Range("D9").Value2 = "0.552"
$ws.Range("E9").Value2 = "  +3.52%  "

$ws.Range("D10").Value2 = "40.55"
$ws.Range("E10").Value2 = "  +5.30%  "

$ws.Range("D11").Value2 = "20.42"
$ws.Range("E11").Value2 = "  +12.27%  "

$ws.Range("E12").Value2 = "  +1.72%  "

$ws.Range("E13").Value2 = "  +1.17%  "

$ws.Range("E14").Value2 = "  +2.68%  "

$ws.Range("D15").Value2 = "2.924.89"
$ws.Range("E15").Value2 = "  +1.69%  "

$ws.Range("D16").Value2 = "2.527.13"
$ws.Range("E16").Value2 = "  +1.71%  "

$ws.Range("E17").Value2 = "  +1.32%  "

$ws.Range("D18").Value2 = "48.230.80"
$ws.Range("E18").Value2 = "  +2.33%  "

$ws.Range("D19").Value2 = "13.40"
$ws.Range("E19").Value2 = "  +5.23%  "

$ws.Range("D20").Value2 = "6.62"
$ws.Range("E20").Value2 = "  +0.30%  "

$ws.Range("E21").Value2 = "  +1.86%  "

$ws.Range("D22").Value2 = "2.69"
$ws.Range("E22").Value2 = "  -0.77%  "

$ws.Range("D23").Value2 = "71.94"
$ws.Range("E23").Value2 = "  +2.37%  "

$ws.Range("D24").Value2 = "269.10"
$ws.Range("E24").Value2 = "  +9.76%  "

$ws.Range("E25").Value2 = "  +0.64%  "

$ws.Range("E26").Value2 = "  +0.05%  "

$ws.Range("D27").Value2 = "26.05"
$ws.Range("E27").Value2 = "  +1.58%  "

$ws.Range("B28").Value2 = "Cosmos"
$ws.Range("C28").Value2 = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D28").Value2 = "10.14"
$ws.Range("E28").Value2 = "  +1.47%  "

$ws.Range("B29").Value2 = "Toncoin"
$ws.Range("C29").Value2 = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value2 = "2.25"
$ws.Range("E29").Value2 = "  -1.15%  "

$ws.Range("E30").Value2 = "  +7.20%  "

$ws.Range("D31").Value2 = "35.69"
$ws.Range("E31").Value2 = "  +3.86%  "

$ws.Range("E32").Value2 = "  +0.49%  "

$ws.Range("D33").Value2 = "19.71"
$ws.Range("E33").Value2 = "  -2.30%  "

$ws.Range("E34").Value2 = "  +1.24%  "

$ws.Range("D35").Value2 = "1.01"
$ws.Range("E35").Value2 = "  +0.10%  "

$ws.Range("D36").Value2 = "0.0787"
$ws.Range("E36").Value2 = "  +1.20%  "

$ws.Range("D37").Value2 = "1.99"
$ws.Range("E37").Value2 = "  +1.82%  "

$ws.Range("E38").Value2 = "  +2.18%  "

$ws.Range("D39").Value2 = "3.03"
$ws.Range("E39").Value2 = "  +3.19%  "

$ws.Range("E40").Value2 = "  +0.82%  "

$ws.Range("D41").Value2 = "121.81"
$ws.Range("E41").Value2 = "  +2.22%  "

$ws.Range("B42").Value2 = "EnergySwap"
$ws.Range("C42").Value2 = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D42").Value2 = "22.02"
$ws.Range("E42").Value2 = "  -1.38%  "

$ws.Range("B43").Value2 = "WEMIXToken"
$ws.Range("C43").Value2 = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D43").Value2 = "2.20"
$ws.Range("E43").Value2 = "  -0.77%  "

$ws.Range("E44").Value2 = "  +2.41%  "

$ws.Range("D45").Value2 = "2.027.24"
$ws.Range("E45").Value2 = "  +2.28%  "

$ws.Range("D46").Value2 = "3.16"
$ws.Range("E46").Value2 = "  +5.40%  "

$ws.Range("D47").Value2 = "1.90"
$ws.Range("E47").Value2 = "  +8.16%  "

$ws.Range("E48").Value2 = "  +3.19%  "

$ws.Range("E49").Value2 = "  +1.08%  "

$ws.Range("D50").Value2 = "5.23"
$ws.Range("E50").Value2 = "  +2.51%  "

$ws.Range("D51").Value2 = "79.35"
$ws.Range("E51").Value2 = "  +3.53%  "
